$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles/number formats/borders) of the last existing
# data row (357) down onto the new rows (358-366) before writing values.
$ws.Range("A357:D357").Copy()
$ws.Range("A358:D366").PasteSpecial(-4122)

# New daily data through 1/09/2021 (date serials 44432-44440)
$data = @(
    @(44432, 3, 61, 84.29023476903093),
    @(44433, 1, 54, 74.61758487750279),
    @(44434, 18, 61, 84.29023476903093),
    @(44435, 11, 67, 92.58107753319791),
    @(44436, 21, 72, 99.49011317000372),
    @(44437, 4, 70, 96.7264989152814),
    @(44438, 49, 107, 147.8533626276444),
    @(44439, 9, 113, 156.1442053918114),
    @(44440, 2, 114, 157.5260125191726)
)

$r = 358
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
